$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Junio de 2020 a las 14:12"

# Refreshed country data, re-sorted descending by "Casos totales" (col B).
# Some rows keep the same country with new numbers; others show a new country
# that overtook / fell behind a neighbour in the ranking, shifting the column-A
# label for that row while the row number (ranking position) stays fixed.
$updates = @(
    @{ Row = 4; Country = "Estados Unidos"; Vals = @(2162406, 178, 870077, 1174470, 0, 1, 117859) }
    @{ Row = 7; Country = "India"; Vals = @(333380, 597, 169987, 153869, 0, 4, 9524) }
    @{ Row = 13; Country = "Alemania"; Vals = @(187706, 35, 172600, 6233, 0, 3, 8873) }
    @{ Row = 26; Country = "Bielorrusia"; Vals = @(54680, 707, 30420, 23948, 0, 4, 312) }
    @{ Row = 27; Country = "Suecia"; Vals = @(52383, 139, 0, 0, 0, 17, 4891) }
    @{ Row = 29; Country = "Paises Bajos"; Vals = @(48948, 165, 0, 0, 0, 6, 6065) }
    @{ Row = 56; Country = "Kazajistan"; Vals = @(14809, 313, 9241, 5491, 0, 0, 77) }
    @{ Row = 58; Country = "Dinamarca"; Vals = @(12217, 24, 11090, 529, 0, 1, 598) }
    @{ Row = 73; Country = "Finlandia"; Vals = @(7108, 4, 6200, 582, 0, 0, 326) }
    @{ Row = 76; Country = "Uzbekistan"; Vals = @(5103, 23, 3996, 1088, 0, 0, 19) }
    @{ Row = 83; Country = "Republica de Macedonia"; Vals = @(4157, 100, 1723, 2241, 0, 5, 193) }
    @{ Row = 84; Country = "Hungria"; Vals = @(4076, 7, 2485, 1028, 0, 1, 563) }
    @{ Row = 85; Country = "Luxemburgo"; Vals = @(4070, 0, 3929, 31, 0, 0, 110) }
    @{ Row = 98; Country = "Croacia"; Vals = @(2254, 2, 2140, 7, 0, 0, 107) }
    @{ Row = 114; Country = "Libano"; Vals = @(1464, 18, 875, 557, 0, 0, 32) }
    @{ Row = 115; Country = "Nicaragua"; Vals = @(1464, 0, 953, 456, 0, 0, 55) }
    @{ Row = 116; Country = "Guinea-Bisau"; Vals = @(1460, 0, 153, 1292, 0, 0, 15) }
    @{ Row = 139; Country = "San Marino"; Vals = @(694, 0, 575, 77, 0, 0, 42) }
    @{ Row = 141; Country = "Malta"; Vals = @(650, 1, 603, 38, 0, 0, 9) }
    @{ Row = 151; Country = "Benin"; Vals = @(483, 41, 232, 242, 0, 3, 9) }
    @{ Row = 206; Country = "Islas Malvinas"; Vals = @(13, 0, 13, 0, 0, 0, 0) }
    @{ Row = 207; Country = "Groenlandia"; Vals = @(13, 0, 13, 0, 0, 0, 0) }
    @{ Row = 208; Country = "Santa Sede"; Vals = @(12, 0, 12, 0, 0, 0, 0) }
    @{ Row = 209; Country = "Islas Turcas y Caicos"; Vals = @(12, 0, 11, 0, 0, 0, 1) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Country
    for ($i = 0; $i -lt $u.Vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $u.Vals[$i]
    }
}
